$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '293.03'
$ws.Cells.Item(2, 4).ClearFormats()
$ws.Cells.Item(2, 5).NumberFormat = '@'
$ws.Cells.Item(2, 5).Value = '0.41%'
$ws.Cells.Item(2, 5).ClearFormats()
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '40.41'
$ws.Cells.Item(3, 4).ClearFormats()
$ws.Cells.Item(3, 5).NumberFormat = '@'
$ws.Cells.Item(3, 5).Value = '1.16%'
$ws.Cells.Item(3, 5).ClearFormats()
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '5.006'
$ws.Cells.Item(4, 4).ClearFormats()
$ws.Cells.Item(4, 5).NumberFormat = '@'
$ws.Cells.Item(4, 5).Value = '-0.31%'
$ws.Cells.Item(4, 5).ClearFormats()
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '0.07354'
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).NumberFormat = '@'
$ws.Cells.Item(5, 5).Value = '-0.47%'
$ws.Cells.Item(5, 5).ClearFormats()
$ws.Cells.Item(6, 2).NumberFormat = '@'
$ws.Cells.Item(6, 2).Value = 'GateToken'
$ws.Cells.Item(6, 2).ClearFormats()
$ws.Cells.Item(6, 3).NumberFormat = '@'
$ws.Cells.Item(6, 3).Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Cells.Item(6, 3).ClearFormats()
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '4.288'
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).NumberFormat = '@'
$ws.Cells.Item(6, 5).Value = '-0.85%'
$ws.Cells.Item(6, 5).ClearFormats()
$ws.Cells.Item(7, 2).NumberFormat = '@'
$ws.Cells.Item(7, 2).Value = 'FTXToken'
$ws.Cells.Item(7, 2).ClearFormats()
$ws.Cells.Item(7, 3).NumberFormat = '@'
$ws.Cells.Item(7, 3).Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Cells.Item(7, 3).ClearFormats()
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '1.563'
$ws.Cells.Item(7, 4).ClearFormats()
$ws.Cells.Item(7, 5).NumberFormat = '@'
$ws.Cells.Item(7, 5).Value = '3.71%'
$ws.Cells.Item(7, 5).ClearFormats()
$ws.Cells.Item(8, 2).NumberFormat = '@'
$ws.Cells.Item(8, 2).Value = 'MXToken'
$ws.Cells.Item(8, 2).ClearFormats()
$ws.Cells.Item(8, 3).NumberFormat = '@'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(8, 3).ClearFormats()
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.9243'
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(8, 5).NumberFormat = '@'
$ws.Cells.Item(8, 5).Value = '0.07%'
$ws.Cells.Item(8, 5).ClearFormats()
$ws.Cells.Item(9, 2).NumberFormat = '@'
$ws.Cells.Item(9, 2).Value = 'BTSEToken'
$ws.Cells.Item(9, 2).ClearFormats()
$ws.Cells.Item(9, 3).NumberFormat = '@'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Cells.Item(9, 3).ClearFormats()
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '2.352'
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(9, 5).NumberFormat = '@'
$ws.Cells.Item(9, 5).Value = '-1.96%'
$ws.Cells.Item(9, 5).ClearFormats()
$ws.Cells.Item(10, 2).NumberFormat = '@'
$ws.Cells.Item(10, 2).Value = 'LiechtensteinCryptoassetsExchange'
$ws.Cells.Item(10, 2).ClearFormats()
$ws.Cells.Item(10, 3).NumberFormat = '@'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Cells.Item(10, 3).ClearFormats()
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.1176'
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).NumberFormat = '@'
$ws.Cells.Item(10, 5).Value = '1.45%'
$ws.Cells.Item(10, 5).ClearFormats()
$ws.Cells.Item(11, 2).NumberFormat = '@'
$ws.Cells.Item(11, 2).Value = 'WazirX'
$ws.Cells.Item(11, 2).ClearFormats()
$ws.Cells.Item(11, 3).NumberFormat = '@'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Cells.Item(11, 3).ClearFormats()
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.1814'
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).NumberFormat = '@'
$ws.Cells.Item(11, 5).Value = '3.94%'
$ws.Cells.Item(11, 5).ClearFormats()
$ws.Cells.Item(12, 2).NumberFormat = '@'
$ws.Cells.Item(12, 2).Value = 'BitrueCoin'
$ws.Cells.Item(12, 2).ClearFormats()
$ws.Cells.Item(12, 3).NumberFormat = '@'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Cells.Item(12, 3).ClearFormats()
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.04381'
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(12, 5).NumberFormat = '@'
$ws.Cells.Item(12, 5).Value = '4.66%'
$ws.Cells.Item(12, 5).ClearFormats()
$ws.Cells.Item(13, 2).NumberFormat = '@'
$ws.Cells.Item(13, 2).Value = 'MandalaExchangeToken'
$ws.Cells.Item(13, 2).ClearFormats()
$ws.Cells.Item(13, 3).NumberFormat = '@'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Cells.Item(13, 3).ClearFormats()
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '0.08830'
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).NumberFormat = '@'
$ws.Cells.Item(13, 5).Value = '2.23%'
$ws.Cells.Item(13, 5).ClearFormats()
$ws.Cells.Item(14, 2).NumberFormat = '@'
$ws.Cells.Item(14, 2).Value = 'BitMartToken'
$ws.Cells.Item(14, 2).ClearFormats()
$ws.Cells.Item(14, 3).NumberFormat = '@'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Cells.Item(14, 3).ClearFormats()
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '0.1056'
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(14, 5).NumberFormat = '@'
$ws.Cells.Item(14, 5).Value = '0.34%'
$ws.Cells.Item(14, 5).ClearFormats()
$ws.Cells.Item(15, 2).NumberFormat = '@'
$ws.Cells.Item(15, 2).Value = 'BitForexToken'
$ws.Cells.Item(15, 2).ClearFormats()
$ws.Cells.Item(15, 3).NumberFormat = '@'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Cells.Item(15, 3).ClearFormats()
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.001268'
$ws.Cells.Item(15, 4).ClearFormats()
$ws.Cells.Item(15, 5).NumberFormat = '@'
$ws.Cells.Item(15, 5).Value = '-1.05%'
$ws.Cells.Item(15, 5).ClearFormats()
$ws.Cells.Item(16, 2).NumberFormat = '@'
$ws.Cells.Item(16, 2).Value = 'TigerCash'
$ws.Cells.Item(16, 2).ClearFormats()
$ws.Cells.Item(16, 3).NumberFormat = '@'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Cells.Item(16, 3).ClearFormats()
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '0.005953'
$ws.Cells.Item(16, 4).ClearFormats()
$ws.Cells.Item(16, 5).NumberFormat = '@'
$ws.Cells.Item(16, 5).Value = '0.16%'
$ws.Cells.Item(16, 5).ClearFormats()
$ws.Cells.Item(17, 2).NumberFormat = '@'
$ws.Cells.Item(17, 2).Value = 'LEO'
$ws.Cells.Item(17, 2).ClearFormats()
$ws.Cells.Item(17, 3).NumberFormat = '@'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(17, 3).ClearFormats()
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '3.350'
$ws.Cells.Item(17, 4).ClearFormats()
$ws.Cells.Item(17, 5).NumberFormat = '@'
$ws.Cells.Item(17, 5).Value = '-0.60%'
$ws.Cells.Item(17, 5).ClearFormats()
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '0.3306'
$ws.Cells.Item(18, 4).ClearFormats()
$ws.Cells.Item(18, 5).NumberFormat = '@'
$ws.Cells.Item(18, 5).Value = '-0.21%'
$ws.Cells.Item(18, 5).ClearFormats()
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '7.816'
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).NumberFormat = '@'
$ws.Cells.Item(19, 5).Value = '3.18%'
$ws.Cells.Item(19, 5).ClearFormats()
$ws.Cells.Item(20, 5).NumberFormat = '@'
$ws.Cells.Item(20, 5).Value = '2.21%'
$ws.Cells.Item(20, 5).ClearFormats()
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '0.03919'
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).NumberFormat = '@'
$ws.Cells.Item(22, 5).Value = '2.27%'
$ws.Cells.Item(22, 5).ClearFormats()
$ws.Cells.Item(23, 5).NumberFormat = '@'
$ws.Cells.Item(23, 5).Value = '-2.01%'
$ws.Cells.Item(23, 5).ClearFormats()
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '0.003677'
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).NumberFormat = '@'
$ws.Cells.Item(24, 5).Value = '-5.65%'
$ws.Cells.Item(24, 5).ClearFormats()
$ws.Cells.Item(25, 5).NumberFormat = '@'
$ws.Cells.Item(25, 5).Value = '-8.19%'
$ws.Cells.Item(25, 5).ClearFormats()
$ws.Cells.Item(26, 5).NumberFormat = '@'
$ws.Cells.Item(26, 5).Value = '-0.61%'
$ws.Cells.Item(26, 5).ClearFormats()
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.02346'
$ws.Cells.Item(38, 4).ClearFormats()
$ws.Cells.Item(38, 5).NumberFormat = '@'
$ws.Cells.Item(38, 5).Value = '1.32%'
$ws.Cells.Item(38, 5).ClearFormats()
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.05106'
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).NumberFormat = '@'
$ws.Cells.Item(39, 5).Value = '2.33%'
$ws.Cells.Item(39, 5).ClearFormats()
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.005966'
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(40, 5).NumberFormat = '@'
$ws.Cells.Item(40, 5).Value = '48.55%'
$ws.Cells.Item(40, 5).ClearFormats()
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.007869'
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(41, 5).NumberFormat = '@'
$ws.Cells.Item(41, 5).Value = '1.97%'
$ws.Cells.Item(41, 5).ClearFormats()
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.1291'
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).NumberFormat = '@'
$ws.Cells.Item(42, 5).Value = '1.37%'
$ws.Cells.Item(42, 5).ClearFormats()
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.007380'
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).NumberFormat = '@'
$ws.Cells.Item(43, 5).Value = '-0.77%'
$ws.Cells.Item(43, 5).ClearFormats()
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.008057'
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).NumberFormat = '@'
$ws.Cells.Item(44, 5).Value = '1.89%'
$ws.Cells.Item(44, 5).ClearFormats()
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.2910'
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).NumberFormat = '@'
$ws.Cells.Item(45, 5).Value = '-7.58%'
$ws.Cells.Item(45, 5).ClearFormats()
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '0.00006234'
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).NumberFormat = '@'
$ws.Cells.Item(46, 5).Value = '-3.95%'
$ws.Cells.Item(46, 5).ClearFormats()
$ws.Cells.Item(47, 5).NumberFormat = '@'
$ws.Cells.Item(47, 5).Value = '-0.54%'
$ws.Cells.Item(47, 5).ClearFormats()
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '0.04750'
$ws.Cells.Item(48, 4).ClearFormats()
$ws.Cells.Item(48, 5).NumberFormat = '@'
$ws.Cells.Item(48, 5).Value = '-81.14%'
$ws.Cells.Item(48, 5).ClearFormats()
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '0.00002101'
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).NumberFormat = '@'
$ws.Cells.Item(49, 5).Value = '-0.54%'
$ws.Cells.Item(49, 5).ClearFormats()
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '0.0002001'
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(50, 5).NumberFormat = '@'
$ws.Cells.Item(50, 5).Value = '-0.54%'
$ws.Cells.Item(50, 5).ClearFormats()
